$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the style of the existing header row (bold, centered, bordered)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill season record values for each data row (2 through 42)
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 82
    $ws.Cells.Item($r, 32).Value = 0
}
